$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARCHITECTURE")

# Update WWR (window-to-wall ratio) values for SG Office (row 5), Retail (row 6)
# and Restaurant (row 8) from 0.35 to 0.59, across the four orientation
# columns: G (wwr_north), H (wwr_south), I (wwr_east), J (wwr_west).
foreach ($row in 5, 6, 8) {
    $rangeAddress = "G" + $row + ":J" + $row
    $ws.Range($rangeAddress).Value = 0.59
}

# Update the selection on the sheet to match the authored state.
$ws.Activate()
$ws.Range("G8").Select()
